$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2")
$rng.NumberFormat = "@"
$rng.Value = "65.650.30"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$rng = $ws.Range("D3")
$rng.NumberFormat = "@"
$rng.Value = "3.412.43"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  +0.11%  "
$rng = $ws.Range("D5")
$rng.NumberFormat = "@"
$rng.Value = "182.83"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E5").Value = "  -10.49%  "
$rng = $ws.Range("D6")
$rng.NumberFormat = "@"
$rng.Value = "531.24"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E6").Value = "  -5.92%  "
$rng = $ws.Range("D7")
$rng.NumberFormat = "@"
$rng.Value = "0.614"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E7").Value = "  -1.25%  "
$rng = $ws.Range("D8")
$rng.NumberFormat = "@"
$rng.Value = "3.400.41"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E8").Value = "  -5.58%  "
$ws.Range("E9").Value = "  -0.01%  "
$rng = $ws.Range("D10")
$rng.NumberFormat = "@"
$rng.Value = "0.629"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E10").Value = "  -6.51%  "
$rng = $ws.Range("D11")
$rng.NumberFormat = "@"
$rng.Value = "57.67"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E11").Value = "  -5.20%  "
$rng = $ws.Range("D12")
$rng.NumberFormat = "@"
$rng.Value = "0.135"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E12").Value = "  -10.59%  "
$rng = $ws.Range("D13")
$rng.NumberFormat = "@"
$rng.Value = "0.0000256"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E13").Value = "  -10.60%  "
$rng = $ws.Range("D14")
$rng.NumberFormat = "@"
$rng.Value = "9.37"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E14").Value = "  -6.42%  "
$rng = $ws.Range("D15")
$rng.NumberFormat = "@"
$rng.Value = "3.962.97"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E15").Value = "  -5.35%  "
$rng = $ws.Range("D17")
$rng.NumberFormat = "@"
$rng.Value = "3.408.63"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E17").Value = "  -5.61%  "
$rng = $ws.Range("D18")
$rng.NumberFormat = "@"
$rng.Value = "65.443.47"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E18").Value = "  -3.63%  "
$rng = $ws.Range("D19")
$rng.NumberFormat = "@"
$rng.Value = "17.65"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E19").Value = "  -6.25%  "
$rng = $ws.Range("D20")
$rng.NumberFormat = "@"
$rng.Value = "11.29"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E20").Value = "  -8.61%  "
$rng = $ws.Range("D21")
$rng.NumberFormat = "@"
$rng.Value = "0.985"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E21").Value = "  -8.29%  "
$rng = $ws.Range("D22")
$rng.NumberFormat = "@"
$rng.Value = "379.27"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E22").Value = "  -5.66%  "
$rng = $ws.Range("D23")
$rng.NumberFormat = "@"
$rng.Value = "83.38"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E23").Value = "  -2.23%  "
$rng = $ws.Range("D24")
$rng.NumberFormat = "@"
$rng.Value = "3.75"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E24").Value = "  -9.63%  "
$rng = $ws.Range("D25")
$rng.NumberFormat = "@"
$rng.Value = "10.86"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E25").Value = "  -16.86%  "
$rng = $ws.Range("D26")
$rng.NumberFormat = "@"
$rng.Value = "11.63"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E26").Value = "  -7.58%  "
$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$rng = $ws.Range("D27")
$rng.NumberFormat = "@"
$rng.Value = "2.67"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E27").Value = "  -8.84%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$rng = $ws.Range("D28")
$rng.NumberFormat = "@"
$rng.Value = "3.60"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E28").Value = "  -9.60%  "
$rng = $ws.Range("D29")
$rng.NumberFormat = "@"
$rng.Value = "8.56"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E29").Value = "  -9.06%  "
$rng = $ws.Range("D30")
$rng.NumberFormat = "@"
$rng.Value = "682.39"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E30").Value = "  +0.94%  "
$rng = $ws.Range("D31")
$rng.NumberFormat = "@"
$rng.Value = "29.94"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E31").Value = "  -5.24%  "
$rng = $ws.Range("D32")
$rng.NumberFormat = "@"
$rng.Value = "6.74"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E32").Value = "  -19.40%  "
$rng = $ws.Range("D33")
$rng.NumberFormat = "@"
$rng.Value = "11.25"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E33").Value = "  -7.88%  "
$rng = $ws.Range("D34")
$rng.NumberFormat = "@"
$rng.Value = "61.55"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E34").Value = "  -3.51%  "
$ws.Range("E35").Value = "  -6.88%  "
$rng = $ws.Range("D36")
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E36").Value = "  -0.14%  "
$rng = $ws.Range("D37")
$rng.NumberFormat = "@"
$rng.Value = "36.83"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E37").Value = "  -12.89%  "
$rng = $ws.Range("D38")
$rng.NumberFormat = "@"
$rng.Value = "0.388"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E38").Value = "  -8.62%  "
$rng = $ws.Range("D39")
$rng.NumberFormat = "@"
$rng.Value = "0.999"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$rng = $ws.Range("D40")
$rng.NumberFormat = "@"
$rng.Value = "0.127"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E40").Value = "  -6.42%  "
$rng = $ws.Range("D41")
$rng.NumberFormat = "@"
$rng.Value = "2.896.90"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E41").Value = "  -11.88%  "
$rng = $ws.Range("D42")
$rng.NumberFormat = "@"
$rng.Value = "2.77"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E42").Value = "  -13.22%  "
$ws.Range("E43").Value = "  -1.79%  "
$rng = $ws.Range("D44")
$rng.NumberFormat = "@"
$rng.Value = "0.0₃0628"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E44").Value = "  -18.09%  "
$rng = $ws.Range("D45")
$rng.NumberFormat = "@"
$rng.Value = "0.0392"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E45").Value = "  -6.32%  "
$rng = $ws.Range("D46")
$rng.NumberFormat = "@"
$rng.Value = "2.34"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E46").Value = "  -15.16%  "
$rng = $ws.Range("D47")
$rng.NumberFormat = "@"
$rng.Value = "0.126"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E47").Value = "  -3.77%  "
$rng = $ws.Range("D48")
$rng.NumberFormat = "@"
$rng.Value = "134.92"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E48").Value = "  -3.54%  "
$ws.Range("E49").Value = "  -7.75%  "
$rng = $ws.Range("D50")
$rng.NumberFormat = "@"
$rng.Value = "2.57"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E50").Value = "  -5.82%  "
$rng = $ws.Range("D51")
$rng.NumberFormat = "@"
$rng.Value = "7.68"
$rng.NumberFormat = "General"
$rng.Style = "Normal"
$ws.Range("E51").Value = "  -12.99%  "
